$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $wb.Styles.Add("Normal 2")
$style.Font.Name = "Calibri"

for ($r = 1; $r -le 18; $r++) {
    $ws.Cells.Item($r, 5).Value = "Setosa"
}

$ws.Range("E1:E18").Style = "Normal 2"

$ws.Range("K5").Select() | Out-Null
